# WebForm User Assignment execution
# Assigns newly-generated phone numbers (PN_Value, column F) to the
# WebForm test rows on Sheet1. Row 5's F cell previously had no value at
# all and now gets one too.
#
# Values are entered with a leading apostrophe so Excel stores them as
# text (matching the existing column, which holds numeric-looking phone
# numbers as strings) and the style is reset to "Normal" immediately
# after so the apostrophe doesn't leave a quote-prefixed number format
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phoneNumbers = @{
    "F2"  = "9840088107"
    "F3"  = "9840087288"
    "F4"  = "9840075406"
    "F5"  = "9840063320"
    "F6"  = "9840031530"
    "F7"  = "9840055858"
    "F8"  = "9840028207"
    "F9"  = "9840033067"
    "F10" = "9840063673"
}

foreach ($cellRef in $phoneNumbers.Keys) {
    $ws.Range($cellRef).Value = "'" + $phoneNumbers[$cellRef]
    $ws.Range($cellRef).Style = "Normal"
}
